# "montagem de torres removida"
# Remove the "Conclusao Torres" (tower completion) outorgado/realizado rows'
# values from column B, and reset the sheet view scroll/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the values in B17 (DatConclusaoTorresOutorgado) and
# B18 (DatConclusaoTorresRealizado) - the "montagem de torres" fields.
$ws.Range("B17").ClearContents()
$ws.Range("B18").ClearContents()

# Reset view: select D8 and scroll back to the top-left (A1) of the sheet.
$ws.Range("D8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
